$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3698.6858
$ws.Range("I28").Value = 3852.5715
$ws.Range("J28").Value = 3083.1428
$ws.Range("K28").Value = 3852.5715
$ws.Range("L28").Value = 3083.1428
$ws.Range("M28").Value = -3367.5715
$ws.Range("N28").Value = -4053.1428
$ws.Range("H51").Value = 4717
$ws.Range("J51").Value = 5155.875
$ws.Range("L51").Value = 5155.875
$ws.Range("N51").Value = -6123.875
$ws.Range("H92").Value = 4403.4707
$ws.Range("I92").Value = 1816.2222
$ws.Range("K92").Value = 1816.2222
$ws.Range("M92").Value = -568.2221999999999
$ws.Range("H112").Value = 1973.8182
$ws.Range("I112").Value = 1299
$ws.Range("J112").Value = 2536.1667
$ws.Range("K112").Value = 3897
$ws.Range("L112").Value = 7608.500100000001
$ws.Range("M112").Value = -2789
$ws.Range("N112").Value = -9824.500100000001
$ws.Range("H138").Value = 3032.8823
$ws.Range("J138").Value = 3541.12
$ws.Range("L138").Value = 10623.36
$ws.Range("N138").Value = -20903.36

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1125.081
$ws.Range("I2").Value = 1125.081
$ws.Range("K2").Value = 1125.081
$ws.Range("M2").Value = -1012.081
$ws.Range("H32").Value = 15888.939
$ws.Range("I32").Value = 6630.722
$ws.Range("K32").Value = 6630.722
$ws.Range("M32").Value = -6343.722
$ws.Range("H43").Value = 90599.39999999999
$ws.Range("I43").Value = 168999
$ws.Range("K43").Value = 168999
$ws.Range("M43").Value = -168686
$ws.Range("H45").Value = 2751.6667
$ws.Range("I45").Value = 4694.3335
$ws.Range("K45").Value = 4694.3335
$ws.Range("M45").Value = -4317.3335
$ws.Range("H61").Value = 3693.1333
$ws.Range("I61").Value = 3414.25
$ws.Range("K61").Value = 3414.25
$ws.Range("M61").Value = -3202.25
$ws.Range("H109").Value = 67165.5
$ws.Range("J109").Value = 67165.5
$ws.Range("L109").Value = 67165.5
$ws.Range("N109").Value = -69939.5
$ws.Range("H116").Value = 1125.081
$ws.Range("I116").Value = 1125.081
$ws.Range("K116").Value = 1125.081
$ws.Range("M116").Value = 1168.919
$ws.Range("H122").Value = 2483.2
$ws.Range("I122").Value = 2407.5386
$ws.Range("K122").Value = 7222.6158
$ws.Range("M122").Value = -4772.6158
$ws.Range("H136").Value = 3693.1333
$ws.Range("I136").Value = 3414.25
$ws.Range("K136").Value = 10242.75
$ws.Range("M136").Value = -7692.75
$ws.Range("H138").Value = 85000
$ws.Range("J138").Value = 85000
$ws.Range("L138").Value = 85000
$ws.Range("N138").Value = -95280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1125.081
$ws.Range("I3").Value = 1125.081
$ws.Range("K3").Value = 1125.081
$ws.Range("M3").Value = -1011.081
$ws.Range("H105").Value = 3381.1667
$ws.Range("I105").Value = 1759.8
$ws.Range("K105").Value = 1759.8
$ws.Range("M105").Value = -12.79999999999995
$ws.Range("H107").Value = 1133
$ws.Range("I107").Value = 1133
$ws.Range("K107").Value = 1133
$ws.Range("M107").Value = 787
$ws.Range("H138").Value = 20000
$ws.Range("I138").Value = 20000
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 20000
$ws.Range("N138").ClearContents()
$ws.Range("M138").Value = -14860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8660.793
$ws.Range("I31").Value = 3856
$ws.Range("K31").Value = 3856
$ws.Range("M31").Value = -3561
$ws.Range("H34").Value = 8660.793
$ws.Range("I34").Value = 3856
$ws.Range("K34").Value = 3856
$ws.Range("M34").Value = -3654
$ws.Range("H86").Value = 4417.6665
$ws.Range("J86").Value = 4683.5
$ws.Range("L86").Value = 4683.5
$ws.Range("N86").Value = -6929.5
$ws.Range("H89").Value = 4417.6665
$ws.Range("J89").Value = 4683.5
$ws.Range("L89").Value = 23417.5
$ws.Range("N89").Value = -34649.5
$ws.Range("H132").Value = 2893.889
$ws.Range("I132").Value = 2635
$ws.Range("K132").Value = 7905
$ws.Range("M132").Value = -5375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 108.666664
$ws.Range("I7").Value = 108.666664
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 325.999992
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H32").Value = 167466.17
$ws.Range("J32").Value = 999
$ws.Range("L32").Value = 2997
$ws.Range("N32").Value = -3563
$ws.Range("H37").Value = 83340.27
$ws.Range("J37").Value = 83340.27
$ws.Range("L37").Value = 250020.81
$ws.Range("N37").Value = -250244.81
$ws.Range("H46").Value = 652537.5
$ws.Range("I46").Value = 856287.5
$ws.Range("K46").Value = 2568862.5
$ws.Range("M46").Value = -2568771.5
$ws.Range("H92").Value = 782.9
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 1065.8
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 3197.4
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -5693.4
$ws.Range("H132").Value = 1390
$ws.Range("I132").Value = 1323
$ws.Range("J132").Value = 1423.5
$ws.Range("K132").Value = 11907
$ws.Range("L132").Value = 12811.5
$ws.Range("M132").Value = -9377
$ws.Range("N132").Value = -17871.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3610.8
$ws.Range("I126").Value = 2275.2
$ws.Range("K126").Value = 6825.599999999999
$ws.Range("M126").Value = -4355.599999999999
$ws.Range("H141").Value = 119724.664
$ws.Range("J141").Value = 39000
$ws.Range("L141").Value = 39000
$ws.Range("N141").Value = -49360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4886.4443
$ws.Range("I7").Value = 4934.8125
$ws.Range("K7").Value = 4934.8125
$ws.Range("M7").Value = -4822.8125
$ws.Range("H22").Value = 1555.8
$ws.Range("I22").Value = 1112.25
$ws.Range("K22").Value = 1112.25
$ws.Range("M22").Value = -817.25
$ws.Range("H27").Value = 1555.8
$ws.Range("I27").Value = 1112.25
$ws.Range("K27").Value = 1112.25
$ws.Range("M27").Value = -1005.25
$ws.Range("H43").Value = 1855409.6
$ws.Range("I43").Value = 40000.918
$ws.Range("J43").Value = 4033900
$ws.Range("K43").Value = 40000.918
$ws.Range("L43").Value = 4033900
$ws.Range("M43").Value = -39807.918
$ws.Range("N43").Value = -4034286
$ws.Range("H46").Value = 1439.56
$ws.Range("J46").Value = 1288.9
$ws.Range("L46").Value = 1288.9
$ws.Range("N46").Value = -1664.9
$ws.Range("H55").Value = 425.83334
$ws.Range("I55").Value = 313.75
$ws.Range("J55").Value = 650
$ws.Range("K55").Value = 313.75
$ws.Range("L55").Value = 650
$ws.Range("M55").Value = -140.75
$ws.Range("N55").Value = -996
$ws.Range("H100").Value = 3187.25
$ws.Range("I100").Value = 3299.6
$ws.Range("K100").Value = 3299.6
$ws.Range("M100").Value = -2758.6
$ws.Range("H122").Value = 3598.111
$ws.Range("I122").Value = 3551.1333
$ws.Range("J122").Value = 3833
$ws.Range("K122").Value = 10653.3999
$ws.Range("L122").Value = 11499
$ws.Range("M122").Value = -8203.3999
$ws.Range("N122").Value = -16399
$ws.Range("H126").Value = 4886.4443
$ws.Range("I126").Value = 4934.8125
$ws.Range("K126").Value = 14804.4375
$ws.Range("M126").Value = -12334.4375
$ws.Range("H136").Value = 2872.5
$ws.Range("I136").Value = 2706.5625
$ws.Range("K136").Value = 8119.6875
$ws.Range("M136").Value = -5569.6875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 9999
$ws.Range("J21").Value = 9999
$ws.Range("L21").Value = 9999
$ws.Range("N21").Value = -10469
$ws.Range("H35").Value = 9999
$ws.Range("J35").Value = 9999
$ws.Range("L35").Value = 9999
$ws.Range("N35").Value = -10579
$ws.Range("H37").Value = 56000
$ws.Range("I37").Value = 50000
$ws.Range("K37").Value = 50000
$ws.Range("M37").Value = -49797
